$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the header (pushes the existing data
# rows down by 2, old row 2 becomes row 4, etc.)
$ws.Rows.Item(2).Resize(2).Insert()
$ws.Range("A2:C3").ClearFormats()

# Values for the two newly inserted rows
$ws.Cells.Item(2,1).Value = -0.07984677950541165
$ws.Cells.Item(2,2).Value = -0.1310317392150563
$ws.Cells.Item(2,3).Value = 0.2751492535074553

$ws.Cells.Item(3,1).Value = -0.09086012840270966
$ws.Cells.Item(3,2).Value = -0.02373796701431304
$ws.Cells.Item(3,3).Value = 0.2084361910820008

# Append 8 new rows of data after the existing data (now ending at row 23)
$ws.Cells.Item(24,1).Value = -2.15471959114075
$ws.Cells.Item(24,2).Value = 0.2546487897634502
$ws.Cells.Item(24,3).Value = -0.9123589172959333

$ws.Cells.Item(25,1).Value = 0.6330445607503264
$ws.Cells.Item(25,2).Value = 0.5185447335243235
$ws.Cells.Item(25,3).Value = 0.3759170770645095

$ws.Cells.Item(26,1).Value = -0.2923502524693779
$ws.Cells.Item(26,2).Value = -0.1014364187916107
$ws.Cells.Item(26,3).Value = 1.366777941584587

$ws.Cells.Item(27,1).Value = -1.484014511108395
$ws.Cells.Item(27,2).Value = -0.6635967791080455
$ws.Cells.Item(27,3).Value = 1.253775984048842

$ws.Cells.Item(28,1).Value = 0.2468122641245533
$ws.Cells.Item(28,2).Value = 0.4217223922411623
$ws.Cells.Item(28,3).Value = 0.5185143599907569

$ws.Cells.Item(29,1).Value = -0.3225175539652509
$ws.Cells.Item(29,2).Value = -0.6644023060798641
$ws.Cells.Item(29,3).Value = -0.07338536779085762

$ws.Cells.Item(30,1).Value = -0.4933383464813202
$ws.Cells.Item(30,2).Value = -0.1615586355328552
$ws.Cells.Item(30,3).Value = 0.263843480497599

$ws.Cells.Item(31,1).Value = 0.4722494284311972
$ws.Cells.Item(31,2).Value = 0.08936246732870884
$ws.Cells.Item(31,3).Value = 0.3762639736135807
